$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: was "ATC05_PagarCuenta" (empresa/id_cuenta/cuenta in D:F) -> becomes the
# new "ATC05_descargarCartola6meses" test which instead carries rut/password (B:C),
# just like row 2 does.
$ws.Range("A6").Value = "ATC05_descargarCartola6meses"
$ws.Range("B6").Value = "175553878"
$ws.Range("C6").Value = "Rojas651"
$ws.Range("D6:F6").Clear()

# Row 7: was just "ATC06" -> becomes the brand new test page
# "ATC06_descargarCartolaLuz", also using rut/password.
$ws.Range("A7").Value = "ATC06_descargarCartolaLuz"
$ws.Range("B7").Value = "175553878"
$ws.Range("C7").Value = "Rojas651"

# Rows 8/9 keep their previous text ("ATC07"/"ATC08") - nothing else to change there.

# Column A is now wider and manually sized (no longer relying on AutoFit/BestFit).
$ws.Columns.Item(1).ColumnWidth = 28.307291666666668

# Update the active selection to D5.
$ws.Range("D5").Select()
